$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 159528
$ws.Range("C4").Value = 150575
$ws.Range("C7").Value = 5.61
$ws.Range("C8").Value = 64.34999999999999
